$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $r = $ws.Range($CellRef)
    # Prefix with an apostrophe so Excel doesn't coerce number-like
    # strings (e.g. "1.00", "3.00") into numeric values, then strip
    # any resulting formatting tweak (quotePrefix) so the cell keeps
    # its original (default) style - matches source data which is all
    # plain text / inline strings with no explicit style override.
    $r.Value = "'" + $NewValue
    $r.ClearFormats()
}

Set-TextValue "D2" "41.275.12"
Set-TextValue "E2" "  -3.22%  "
Set-TextValue "D3" "2.461.09"
Set-TextValue "E3" "  -2.57%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "309.71"
Set-TextValue "E5" "  +0.66%  "
Set-TextValue "D6" "93.38"
Set-TextValue "E6" "  -6.75%  "
Set-TextValue "D7" "0.551"
Set-TextValue "E7" "  -2.89%  "
Set-TextValue "E8" "  +0.11%  "
Set-TextValue "D9" "0.497"
Set-TextValue "E9" "  -4.92%  "
Set-TextValue "D10" "33.25"
Set-TextValue "E10" "  -6.69%  "
Set-TextValue "D11" "0.0776"
Set-TextValue "E11" "  -3.56%  "
Set-TextValue "E12" "  -0.82%  "
Set-TextValue "D13" "6.97"
Set-TextValue "E13" "  -4.46%  "
Set-TextValue "D14" "2.841.16"
Set-TextValue "E14" "  -2.44%  "
Set-TextValue "D15" "2.488.23"
Set-TextValue "E15" "  -2.07%  "
Set-TextValue "D16" "14.65"
Set-TextValue "E16" "  -5.31%  "
Set-TextValue "D17" "0.780"
Set-TextValue "E17" "  -3.12%  "
Set-TextValue "D18" "41.275.33"
Set-TextValue "E18" "  -3.12%  "
Set-TextValue "D19" "6.30"
Set-TextValue "E19" "  -6.35%  "
Set-TextValue "D20" "0.0₃0918"
Set-TextValue "E20" "  -3.26%  "
Set-TextValue "D21" "11.31"
Set-TextValue "E21" "  -6.99%  "
Set-TextValue "D22" "67.95"
Set-TextValue "E22" "  -1.92%  "
Set-TextValue "D23" "235.90"
Set-TextValue "E23" "  -2.84%  "
Set-TextValue "E24" "  -3.87%  "
Set-TextValue "E25" "  +0.03%  "
Set-TextValue "E26" "  -6.31%  "
Set-TextValue "D27" "23.98"
Set-TextValue "E27" "  -5.72%  "
Set-TextValue "E28" "  -5.60%  "
Set-TextValue "D29" "9.60"
Set-TextValue "E29" "  -5.05%  "
Set-TextValue "D30" "35.62"
Set-TextValue "E30" "  -7.84%  "
Set-TextValue "D31" "151.71"
Set-TextValue "E31" "  -3.79%  "
Set-TextValue "D32" "5.50"
Set-TextValue "E32" "  -4.13%  "
Set-TextValue "E33" "  -5.78%  "
Set-TextValue "D34" "2.57"
Set-TextValue "E34" "  -2.57%  "
Set-TextValue "D35" "0.0737"
Set-TextValue "E35" "  -5.88%  "
Set-TextValue "D36" "3.00"
Set-TextValue "E36" "  -5.60%  "
Set-TextValue "B37" "ARBITRUM"
Set-TextValue "C37" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D37" "1.87"
Set-TextValue "E37" "  -6.23%  "
Set-TextValue "B38" "Celestia"
Set-TextValue "C38" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D38" "16.84"
Set-TextValue "E38" "  -5.27%  "
Set-TextValue "B39" "Kaspa"
Set-TextValue "C39" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D39" "0.103"
Set-TextValue "E39" "  -7.02%  "
Set-TextValue "B40" "Stellar"
Set-TextValue "C40" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D40" "0.113"
Set-TextValue "E40" "  -3.69%  "
Set-TextValue "D41" "4.15"
Set-TextValue "E41" "  -0.48%  "
Set-TextValue "B42" "EnergySwap"
Set-TextValue "C42" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D42" "20.35"
Set-TextValue "E42" "  -6.68%  "
Set-TextValue "B43" "FirstDigitalUSD"
Set-TextValue "C43" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D43" "1.00"
Set-TextValue "E43" "  +0.14%  "
Set-TextValue "D44" "1.978.72"
Set-TextValue "E44" "  -1.42%  "
Set-TextValue "D45" "0.0283"
Set-TextValue "E45" "  -5.72%  "
Set-TextValue "D46" "3.02"
Set-TextValue "E46" "  -7.41%  "
Set-TextValue "D47" "8.61"
Set-TextValue "E47" "  -2.99%  "
Set-TextValue "D48" "70.64"
Set-TextValue "E48" "  -2.06%  "
Set-TextValue "D49" "96.44"
Set-TextValue "E49" "  -4.70%  "
Set-TextValue "D50" "74.12"
Set-TextValue "E50" "  -6.45%  "
Set-TextValue "E51" "  -6.86%  "
